$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.469.08'
$ws.Range('E2').Value = '  +0.98%  '
$ws.Range('D3').Value = '3.497.19'
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('D5').Value = '599.41'
$ws.Range('E5').Value = '  +0.90%  '
$ws.Range('D6').Value = '180.52'
$ws.Range('E6').Value = '  +4.70%  '
$ws.Range('D7').Value = '0.611'
$ws.Range('E7').Value = '  +4.88%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '3.498.63'
$ws.Range('E9').Value = '  +0.16%  '
$ws.Range('D10').Value = '0.139'
$ws.Range('E10').Value = '  +5.65%  '
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('D13').Value = '4.099.57'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '32.34'
$ws.Range('E14').Value = '  +10.59%  '
$ws.Range('E15').Value = '  +0.28%  '
$ws.Range('D16').Value = '67.453.91'
$ws.Range('E16').Value = '  +0.94%  '
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').Value = '3.493.76'
$ws.Range('E18').Value = '  -0.79%  '
$ws.Range('E19').Value = '  +0.79%  '
$ws.Range('D20').Value = '14.30'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('D21').Value = '390.77'
$ws.Range('E21').Value = '  -0.57%  '
$ws.Range('D22').Value = '7.96'
$ws.Range('E22').Value = '  +0.50%  '
$ws.Range('D23').Value = '73.99'
$ws.Range('E23').Value = '  +1.10%  '
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  +0.06%  '
$ws.Range('E26').Value = '  +0.68%  '
$ws.Range('E27').Value = '  +1.24%  '
$ws.Range('E28').Value = '  +1.59%  '
$ws.Range('E29').Value = '  -2.69%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').Value = '2.08'
$ws.Range('E33').Value = '  +1.31%  '
$ws.Range('D34').Value = '23.59'
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('D35').Value = '7.41'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  -0.25%  '
$ws.Range('D38').Value = '162.71'
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('D39').Value = '0.872'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('D40').Value = '2.82'
$ws.Range('E40').Value = '  +12.04%  '
$ws.Range('D41').Value = '1.89'
$ws.Range('E41').Value = '  -0.50%  '
$ws.Range('D42').Value = '6.84'
$ws.Range('E42').Value = '  +0.09%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '2.858.59'
$ws.Range('E43').Value = '  +0.72%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = '4.65'
$ws.Range('E44').Value = '  +0.39%  '
$ws.Range('D45').Value = '26.52'
$ws.Range('E45').Value = '  +2.00%  '
$ws.Range('D46').Value = '26.85'
$ws.Range('E46').Value = '  -1.19%  '
$ws.Range('E47').Value = '  -1.68%  '
$ws.Range('D48').Value = '41.71'
$ws.Range('E48').Value = '  -2.20%  '
$ws.Range('E49').Value = '  -0.08%  '
$ws.Range('D50').Value = '333.73'
$ws.Range('E50').Value = '  -1.16%  '
$ws.Range('E51').Value = '  -1.00%  '
